# ---------------------------------------------------------------------------
# Applies the cell-level changes described by the upstream OOXML diff:
#   - Sheet 展览   (Exhibition):   F-column (想去人数) bumps, 2 name edits in C12/C14
#   - Sheet 演出   (Performance):  F-column bumps, G10 128 -> "不可售" (sold out/unavailable)
#   - Sheet 本地生活 (Local life):  F-column bumps
#   - Sheet 全部类型 (All types):   F-column bumps, 1 name edit in C15 (mirrors 展览 sheet)
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 68
$ws.Range("F4").Value = 475
$ws.Range("F5").Value = 2014
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 7746
$ws.Range("F8").Value = 246
$ws.Range("F9").Value = 24
$ws.Range("C12").Value = "上海·第五人格同人0nly 2.0"
$ws.Range("F12").Value = 1735
$ws.Range("F13").Value = 1511
$ws.Range("C14").Value = "上海·火影忍者同人only"
$ws.Range("F14").Value = 1304
$ws.Range("F15").Value = 150
$ws.Range("F16").Value = 3804
$ws.Range("F17").Value = 5932
$ws.Range("F18").Value = 664
$ws.Range("F19").Value = 21
$ws.Range("F20").Value = 1051
$ws.Range("F21").Value = 1216
$ws.Range("F22").Value = 398
$ws.Range("F23").Value = 6089
$ws.Range("F26").Value = 4146
$ws.Range("F27").Value = 691
$ws.Range("F28").Value = 1911
$ws.Range("F29").Value = 1147
$ws.Range("F30").Value = 287
$ws.Range("F33").Value = 29
$ws.Range("F34").Value = 191
$ws.Range("F35").Value = 8
$ws.Range("F36").Value = 318
$ws.Range("F37").Value = 1135
$ws.Range("F38").Value = 490
$ws.Range("F39").Value = 1853
$ws.Range("F40").Value = 86
$ws.Range("F41").Value = 389
$ws.Range("F42").Value = 144
$ws.Range("F43").Value = 1101
$ws.Range("F45").Value = 59
$ws.Range("F46").Value = 28
$ws.Range("F48").Value = 164
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 79
$ws.Range("F4").Value = 984
$ws.Range("F5").Value = 119
$ws.Range("F9").Value = 8
$ws.Range("G10").Value = "不可售"
$ws.Range("F11").Value = 661
$ws.Range("F12").Value = 349
$ws.Range("F13").Value = 396
$ws.Range("F15").Value = 199
$ws.Range("F16").Value = 104
$ws.Range("F19").Value = 343
$ws.Range("F20").Value = 158
$ws.Range("F22").Value = 32
$ws.Range("F25").Value = 87
$ws.Range("F29").Value = 267
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1547
$ws.Range("F8").Value = 3067
$ws.Range("F9").Value = 870
$ws.Range("F10").Value = 1030
$ws.Range("F11").Value = 1175
$ws.Range("F12").Value = 1524
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1547
$ws.Range("F4").Value = 475
$ws.Range("F6").Value = 3067
$ws.Range("F7").Value = 2014
$ws.Range("F8").Value = 7746
$ws.Range("F10").Value = 870
$ws.Range("F12").Value = 1735
$ws.Range("F13").Value = 1511
$ws.Range("F14").Value = 1175
$ws.Range("C15").Value = "上海·火影忍者同人only"
$ws.Range("F15").Value = 1304
$ws.Range("F16").Value = 661
$ws.Range("F17").Value = 150
$ws.Range("F18").Value = 1524
$ws.Range("F19").Value = 3804
$ws.Range("F20").Value = 349
$ws.Range("F21").Value = 396
$ws.Range("F22").Value = 664
$ws.Range("F23").Value = 1051
$ws.Range("F24").Value = 1216
$ws.Range("F25").Value = 398
$ws.Range("F26").Value = 6089
$ws.Range("F28").Value = 4146
$ws.Range("F29").Value = 691
$ws.Range("F30").Value = 1911
$ws.Range("F31").Value = 1147
$ws.Range("F32").Value = 287
$ws.Range("F33").Value = 29
$ws.Range("F34").Value = 158
$ws.Range("F35").Value = 191
$ws.Range("F36").Value = 318
$ws.Range("F37").Value = 1135
$ws.Range("F38").Value = 490
$ws.Range("F39").Value = 1853
$ws.Range("F41").Value = 86
$ws.Range("F42").Value = 390
$ws.Range("F43").Value = 1101
$ws.Range("F46").Value = 267
$ws.Range("F48").Value = 164
